# Saldo.xlsx update — "Add files via upload"
#
# The refreshed export drops 11 accounts that no longer appear in the
# source report and adds one new account (TATIANA) right after the
# 005198093 / ANA row. No other rows change position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the rows that are no longer present in the refreshed export.
# Deleted from the bottom up (by original row number) so earlier deletes
# never shift the row number of a delete target still to come.
$ws.Rows(33).Delete()   # 004398174  DANIELE    1538.82
$ws.Rows(31).Delete()   # 004855570  LUISA      1702.58
$ws.Rows(24).Delete()   # 005053939  VIRGILIO   2655.81
$ws.Rows(19).Delete()   # 004641487  LAILA      7541.65
$ws.Rows(18).Delete()   # 004643737  LARA       7574.7
$ws.Rows(14).Delete()   # 004452597  LARA       13287.16
$ws.Rows(11).Delete()   # 004948033  GUILHERME  22000
$ws.Rows(10).Delete()   # 004474776  GILSON     23860.68
$ws.Rows(8).Delete()    # 004643746  MARIO      35422.51
$ws.Rows(4).Delete()    # 004450724  ASSAKO     69863.59
$ws.Rows(2).Delete()    # 004479287  ANA        240492.2

# --- Insert the new account row right after 005198093 / ANA (now row 11).
$ws.Rows(12).Insert()

# Force column A to be stored as text so the leading zeros in the account
# number survive (matches how every other "Conta" cell in the sheet is
# stored).
$ws.Cells.Item(12, 1).Value = "'005366671"
$ws.Cells.Item(12, 2).Value = "TATIANA"
$ws.Cells.Item(12, 3).Value = 6250
